$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:E2").NumberFormat = "@"
$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '61.524.14'
$ws.Range("E2").Value = '  +2.17%  '

$ws.Range("B3:E3").NumberFormat = "@"
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '3.389.18'
$ws.Range("E3").Value = '  +2.87%  '

$ws.Range("B4:E4").NumberFormat = "@"
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("B5:E5").NumberFormat = "@"
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '574.69'
$ws.Range("E5").Value = '  +2.82%  '

$ws.Range("B6:E6").NumberFormat = "@"
$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = '137.49'
$ws.Range("E6").Value = '  +7.45%  '

$ws.Range("B7:E7").NumberFormat = "@"
$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("B8:E8").NumberFormat = "@"
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.388.25'
$ws.Range("E8").Value = '  +2.86%  '

$ws.Range("B9:E9").NumberFormat = "@"
$ws.Range("B9").Value = 'XRP'
$ws.Range("C9").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D9").Value = '0.477'
$ws.Range("E9").Value = '  +0.94%  '

$ws.Range("B10:E10").NumberFormat = "@"
$ws.Range("B10").Value = 'Toncoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D10").Value = '7.46'
$ws.Range("E10").Value = '  +1.84%  '

$ws.Range("B11:E11").NumberFormat = "@"
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.127'
$ws.Range("E11").Value = '  +8.64%  '

$ws.Range("B12:E12").NumberFormat = "@"
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").Value = '0.394'
$ws.Range("E12").Value = '  +6.62%  '

$ws.Range("B13:E13").NumberFormat = "@"
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '3.969.90'
$ws.Range("E13").Value = '  +2.89%  '

$ws.Range("B14:E14").NumberFormat = "@"
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").Value = '0.121'
$ws.Range("E14").Value = '  +1.73%  '

$ws.Range("B15:E15").NumberFormat = "@"
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '0.0000179'
$ws.Range("E15").Value = '  +7.46%  '

$ws.Range("B16:E16").NumberFormat = "@"
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.389.59'
$ws.Range("E16").Value = '  +2.94%  '

$ws.Range("B17:E17").NumberFormat = "@"
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").Value = '25.30'
$ws.Range("E17").Value = '  +4.92%  '

$ws.Range("B18:E18").NumberFormat = "@"
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '61.637.51'
$ws.Range("E18").Value = '  +2.01%  '

$ws.Range("B19:E19").NumberFormat = "@"
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '14.06'
$ws.Range("E19").Value = '  +5.64%  '

$ws.Range("B20:E20").NumberFormat = "@"
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '5.89'
$ws.Range("E20").Value = '  +4.74%  '

$ws.Range("B21:E21").NumberFormat = "@"
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '9.39'
$ws.Range("E21").Value = '  +4.08%  '

$ws.Range("B22:E22").NumberFormat = "@"
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '388.07'
$ws.Range("E22").Value = '  +10.63%  '

$ws.Range("B23:E23").NumberFormat = "@"
$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D23").Value = '0.571'
$ws.Range("E23").Value = '  +3.54%  '

$ws.Range("B24:E24").NumberFormat = "@"
$ws.Range("B24").Value = 'WrappedeETH'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D24").Value = '3.526.01'
$ws.Range("E24").Value = '  +2.98%  '

$ws.Range("B25:E25").NumberFormat = "@"
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.20%  '

$ws.Range("B26:E26").NumberFormat = "@"
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").Value = '0.0000127'
$ws.Range("E26").Value = '  +17.73%  '

$ws.Range("B27:E27").NumberFormat = "@"
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").Value = '71.09'
$ws.Range("E27").Value = '  +2.68%  '

$ws.Range("B28:E28").NumberFormat = "@"
$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").Value = '1.59'
$ws.Range("E28").Value = '  +11.57%  '

$ws.Range("B29:E29").NumberFormat = "@"
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '7.67'
$ws.Range("E29").Value = '  +6.83%  '

$ws.Range("B30:E30").NumberFormat = "@"
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.03%  '

$ws.Range("B31:E31").NumberFormat = "@"
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '8.30'
$ws.Range("E31").Value = '  +6.59%  '

$ws.Range("B32:E32").NumberFormat = "@"
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").Value = '0.158'
$ws.Range("E32").Value = '  +5.57%  '

$ws.Range("B33:E33").NumberFormat = "@"
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = '2.15'
$ws.Range("E33").Value = '  +3.11%  '

$ws.Range("B34:E34").NumberFormat = "@"
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("B35:E35").NumberFormat = "@"
$ws.Range("B35").Value = 'RenzoRestakedETH'
$ws.Range("C35").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D35").Value = '3.419.43'
$ws.Range("E35").Value = '  +2.85%  '

$ws.Range("B36:E36").NumberFormat = "@"
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '23.47'
$ws.Range("E36").Value = '  +3.54%  '

$ws.Range("B37:E37").NumberFormat = "@"
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = '5.49'
$ws.Range("E37").Value = '  +6.05%  '

$ws.Range("B38:E38").NumberFormat = "@"
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").Value = '6.96'
$ws.Range("E38").Value = '  +3.46%  '

$ws.Range("B39:E39").NumberFormat = "@"
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '1.54'
$ws.Range("E39").Value = '  +4.69%  '

$ws.Range("B40:E40").NumberFormat = "@"
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '161.90'
$ws.Range("E40").Value = '  +2.17%  '

$ws.Range("B41:E41").NumberFormat = "@"
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = '0.0796'
$ws.Range("E41").Value = '  +6.94%  '

$ws.Range("B42:E42").NumberFormat = "@"
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.02%  '

$ws.Range("B43:E43").NumberFormat = "@"
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '1.73'
$ws.Range("E43").Value = '  +12.63%  '

$ws.Range("B44:E44").NumberFormat = "@"
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '0.772'
$ws.Range("E44").Value = '  +4.39%  '

$ws.Range("B45:E45").NumberFormat = "@"
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").Value = '4.45'
$ws.Range("E45").Value = '  +2.60%  '

$ws.Range("B46:E46").NumberFormat = "@"
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").Value = '1.22'
$ws.Range("E46").Value = '  +6.58%  '

$ws.Range("B47:E47").NumberFormat = "@"
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").Value = '41.32'
$ws.Range("E47").Value = '  +0.78%  '

$ws.Range("B48:E48").NumberFormat = "@"
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '24.77'
$ws.Range("E48").Value = '  +8.46%  '

$ws.Range("B49:E49").NumberFormat = "@"
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = '6.97'
$ws.Range("E49").Value = '  +4.89%  '

$ws.Range("B50:E50").NumberFormat = "@"
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '22.93'
$ws.Range("E50").Value = '  +7.58%  '

$ws.Range("B51:E51").NumberFormat = "@"
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '2.370.80'
$ws.Range("E51").Value = '  +9.47%  '

